$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value    = "Ready for handoff"

# --- Timestamps bumped forward (new handoff generated) ---
$overview.Range("G2").Value = "2016-08-25 11:02:19"
$zhcn.Range("H2").Value     = "2016-08-25 11:02:13"
$dede.Range("H2").Value     = "2016-08-25 11:02:19"

# --- Column widths shrink from ~30 chars to ~17 chars on the
#     "Status"-ish columns now that the text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
